$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-PlainCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = $val
}

# Row 2
Set-PlainCell 2 4 '26.525.63'
Set-PlainCell 2 5 '  +2.17%  '

# Row 3
Set-PlainCell 3 4 '1.682.74'
Set-PlainCell 3 5 '  +2.69%  '

# Row 4
Set-PlainCell 4 5 '  -0.17%  '

# Row 5
Set-TextCell 5 4 '217.59'
Set-PlainCell 5 5 '  +3.93%  '

# Row 6
Set-TextCell 6 4 '0.5331'
Set-PlainCell 6 5 '  +3.17%  '

# Row 7
Set-TextCell 7 4 '1.001'
Set-PlainCell 7 5 '  -0.17%  '

# Row 8
Set-TextCell 8 4 '0.2678'
Set-PlainCell 8 5 '  +4.69%  '

# Row 9
Set-TextCell 9 4 '0.06425'
Set-PlainCell 9 5 '  +3.06%  '

# Row 10
Set-TextCell 10 4 '21.53'
Set-PlainCell 10 5 '  +5.94%  '

# Row 11
Set-TextCell 11 4 '0.07808'
Set-PlainCell 11 5 '  +3.59%  '

# Row 12
Set-PlainCell 12 2 'Polkadot'
Set-PlainCell 12 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 12 4 '4.510'
Set-PlainCell 12 5 '  +3.53%  '

# Row 13
Set-PlainCell 13 2 'WrappedEther'
Set-PlainCell 13 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-PlainCell 13 4 '1.670.50'
Set-PlainCell 13 5 '  +1.77%  '

# Row 14
Set-TextCell 14 4 '0.5626'
Set-PlainCell 14 5 '  +4.31%  '

# Row 15
Set-PlainCell 15 4 '0.0₅8416'
Set-PlainCell 15 5 '  +6.18%  '

# Row 16
Set-TextCell 16 4 '66.06'
Set-PlainCell 16 5 '  +1.93%  '

# Row 17
Set-PlainCell 17 4 '26.560.59'
Set-PlainCell 17 5 '  +2.22%  '

# Row 18
Set-TextCell 18 4 '1.001'
Set-PlainCell 18 5 '  -0.14%  '

# Row 19
Set-TextCell 19 4 '4.800'
Set-PlainCell 19 5 '  +3.66%  '

# Row 20
Set-TextCell 20 4 '195.93'
Set-PlainCell 20 5 '  +5.95%  '

# Row 21
Set-TextCell 21 4 '10.43'
Set-PlainCell 21 5 '  +4.41%  '

# Row 22
Set-TextCell 22 4 '6.379'
Set-PlainCell 22 5 '  +4.89%  '

# Row 23
Set-PlainCell 23 5 '  -0.18%  '

# Row 24
Set-TextCell 24 4 '143.25'
Set-PlainCell 24 5 '  -1.40%  '

# Row 25
Set-TextCell 25 4 '0.1284'
Set-PlainCell 25 5 '  +8.04%  '

# Row 26
Set-TextCell 26 4 '7.471'
Set-PlainCell 26 5 '  +2.12%  '

# Row 28
Set-PlainCell 28 5 '  +3.10%  '

# Row 29
Set-TextCell 29 4 '0.06126'
Set-PlainCell 29 5 '  +3.02%  '

# Row 30
Set-TextCell 30 4 '1.278'
Set-PlainCell 30 5 '  +2.86%  '

# Row 31
Set-TextCell 31 4 '3.614'
Set-PlainCell 31 5 '  +7.94%  '

# Row 32
Set-TextCell 32 4 '3.460'
Set-PlainCell 32 5 '  +3.57%  '

# Row 33
Set-TextCell 33 4 '1.708'
Set-PlainCell 33 5 '  +6.70%  '

# Row 34
Set-TextCell 34 4 '1.017'
Set-PlainCell 34 5 '  +5.02%  '

# Row 35
Set-PlainCell 35 5 '  +1.49%  '

# Row 36
Set-TextCell 36 4 '2.788'
Set-PlainCell 36 5 '  +2.08%  '

# Row 37
Set-TextCell 37 4 '0.5714'
Set-PlainCell 37 5 '  -2.25%  '

# Row 38
Set-PlainCell 38 5 '  +3.00%  '

# Row 39
Set-TextCell 39 4 '5.956'
Set-PlainCell 39 5 '  +4.54%  '

# Row 40
Set-TextCell 40 4 '0.8732'
Set-PlainCell 40 5 '  +3.83%  '

# Row 41
Set-PlainCell 41 4 '1.059.55'
Set-PlainCell 41 5 '  +1.81%  '

# Row 43
Set-TextCell 43 4 '99.98'
Set-PlainCell 43 5 '  +0.35%  '

# Row 44
Set-PlainCell 44 4 '1.832.80'
Set-PlainCell 44 5 '  +2.22%  '

# Row 45
Set-PlainCell 45 2 'Aave'
Set-PlainCell 45 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 45 4 '57.28'
Set-PlainCell 45 5 '  +5.46%  '

# Row 46
Set-PlainCell 46 2 'EnergySwap'
Set-PlainCell 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 46 4 '8.156'
Set-PlainCell 46 5 '  +2.50%  '

# Row 47
Set-PlainCell 47 2 'Frax'
Set-PlainCell 47 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 47 4 '0.9986'
Set-PlainCell 47 5 '  +0.27%  '

# Row 48
Set-PlainCell 48 2 'Cronos'
Set-PlainCell 48 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 48 4 '0.05202'
Set-PlainCell 48 5 '  +0.25%  '

# Row 49
Set-PlainCell 49 2 'Aptos'
Set-PlainCell 49 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 49 4 '6.086'
Set-PlainCell 49 5 '  +5.45%  '

# Row 50
Set-PlainCell 50 2 'Mantle'
Set-PlainCell 50 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 50 4 '0.4240'
Set-PlainCell 50 5 '  +0.14%  '

# Row 51
Set-PlainCell 51 2 'Algorand'
Set-PlainCell 51 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 51 4 '0.09909'
Set-PlainCell 51 5 '  +3.66%  '
